# Update the "built on ..." timestamp embedded in the version string
# from "January 30 2026 16.19.47 EST" to "February 02 2026 12.49.33 EST"
# across the "About" sheet (A2, A6) and the "Boundaries and methane
# sources" sheet (S2:S32).

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet  = $wb.Worksheets.Item("Boundaries and methane sources")

# --- About sheet: A2 and A6 contain the build timestamp within larger strings ---
$cellA2 = $aboutSheet.Range("A2")
$textA2 = $cellA2.Value()
$cellA2.Value = $textA2.Replace($oldStamp, $newStamp)

$cellA6 = $aboutSheet.Range("A6")
$textA6 = $cellA6.Value()
$cellA6.Value = $textA6.Replace($oldStamp, $newStamp)

# --- Boundaries and methane sources sheet: column S rows 2-32 hold the
#     build_version string, identical to "mines - January 30 (built on
#     <stamp>)" for every data row ---
$usedRange = $dataSheet.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $dataSheet.Cells.Item($r, 19)  # column S = 19
    $val = $cell.Value()
    if ($val -ne $null -and $val.ToString().Contains($oldStamp)) {
        $text = $val.ToString()
        $cell.Value = $text.Replace($oldStamp, $newStamp)
    }
}
